# Conserto do erro com o rótulo da coluna 2050 nas tabelas e retirada das
# linhas com total das tabelas.

$wb = $excel.ActiveWorkbook

# Sheets whose E1 header should read "2050" and whose row 13 ("Total") must
# be removed.
$simpleSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)"
)

foreach ($name in $simpleSheets) {
    $ws = $wb.Worksheets.Item($name)
    # A leading apostrophe keeps Excel from re-interpreting the 4-digit
    # label as a number, so the cell is stored as text (matching the
    # "2015"/"2030"/"2040" labels already next to it).
    $ws.Range("E1").Value = "'2050"
    $ws.Rows.Item(13).Delete()
}

# This sheet uses a "yyyy-yyyy" period label instead of a single year.
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws4.Range("E1").Value = "2041-2050"
$ws4.Rows.Item(13).Delete()

# This sheet already had no "Total" row, only the mislabeled header cell.
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
$ws5.Range("E1").Value = "'2050"

# This sheet has its own "Total" row in a different position (row 4, column B).
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Rows.Item(4).Delete()
